$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,('(USA) Senior, Data Engineer - Senior Software Engineer', 'Sam''s Club', 'Sunnyvale, CA, US USA', 23.3, 'LangChain, RAG, LLaMA, BigQuery, Dataflow, Docker, Kubernetes, CI/CD, Git, BigQuery', '2026-02-24', 'https://www.indeed.com/viewjob?jk=6bedd59c754ec609')
    ,('Sr. MLOps Engineer', 'Hudson Manpower', 'Plano, TX, US USA', 21.1, 'TensorFlow, PyTorch, AWS SageMaker, GCP Vertex AI, Azure ML, MLflow, Docker, Kubernetes, CI/CD, Jenkins', '2026-02-24', 'https://www.indeed.com/viewjob?jk=df841e2c966e93a8')
    ,('Software Engineer Sys 3', 'Lam Research', 'Fremont, CA, US USA', 16.7, 'RAG, Data Lake, FastAPI, Docker, Kubernetes, CI/CD, Jenkins, Git, Kafka, Python', '2026-02-24', 'https://www.indeed.com/viewjob?jk=0b87979ce36653f6')
    ,('Data & ML Engineer', 'EssilorLuxottica', 'Dallas, TX, US USA', 16.7, 'Data Scientist, RAG, Synapse, Data Lake, MLflow, Docker, Kubernetes, AKS, CI/CD, Git', '2026-02-24', 'https://www.indeed.com/viewjob?jk=84539e9d32b8f8db')
    ,('Senior Data Analyst', 'Konrad Group', 'New York, NY, US USA', 16.7, 'RAG, Redshift, BigQuery, Git, Snowflake, Databricks, BigQuery, Redshift, Tableau, Power BI', '2026-02-24', 'https://www.indeed.com/viewjob?jk=85e5ffc39cf3d166')
    ,('Senior Software Engineer, Data Services', 'Fox Corporation', 'Los Angeles, CA, US USA', 14.4, 'Copilot, Kubernetes, AKS, CI/CD, GitHub Actions, Terraform, Git, NoSQL, Python, SQL', '2026-02-24', 'https://www.indeed.com/viewjob?jk=8a2b768ed31bba34')
    ,('Mid-Level AI / Machine Learning Software Engineer', 'MTSI', 'Huntsville, AL, US USA', 14.4, 'RAG, Hugging Face, TensorFlow, PyTorch, Keras, Docker, Git, Matplotlib, Seaborn, Python', '2026-02-23', 'https://www.indeed.com/viewjob?jk=96ede8f4dfa331b6')
    ,('ML Engineer', 'Solventum', 'Austin, PA, US USA', 13.3, 'LangChain, PyTorch, MLflow, Docker, Kubernetes, CI/CD, Git, Python, SQL, R', '2026-02-24', 'https://www.indeed.com/viewjob?jk=a3a0ae84f6b87980')
    ,('Intern the Otsuka Way 2026 - GenAI Engineer Intern', 'Otsuka', 'Princeton, NJ, US USA', 12.2, 'AI Engineer, Generative AI, LangChain, RAG, FAISS, Pinecone, Prompt Engineering, Git, Python, R', '2026-02-24', 'https://www.indeed.com/viewjob?jk=6e718096f3c0913a')
    ,('Senior Analytics Engineer, Sales Analytics (GTM Data Engineering)', 'Klaviyo', 'Boston, MA, US USA', 11.1, 'S3, EC2, Terraform, Git, Snowflake, Tableau, Python, SQL, R, Scala', '2026-02-24', 'https://www.indeed.com/viewjob?jk=721d4b934257a77e')
    ,('Senior Digital Analytics Customer Experience (CX) Engineer', 'Sherwin-Williams', 'Cleveland, OH, US USA', 11.1, 'Git, Snowflake, Databricks, Tableau, Power BI, Python, SQL, R, Scala, A/B Testing', '2026-02-24', 'https://www.indeed.com/viewjob?jk=12fcbbc0dd03080c')
    ,('Senior Data Scientist - Agentic AI', 'Allstate Insurance', 'Remote, US USA', 11.1, 'Data Scientist, RAG, TensorFlow, PyTorch, Docker, CI/CD, Python, SQL, R, Scala', '2026-02-24', 'https://www.indeed.com/viewjob?jk=2ae05a4ee446a6a0')
    ,('Senior Engineer I, Observability', 'DigitalOcean', 'Seattle, WA, US USA', 10, 'RAG, Kubernetes, CI/CD, GitHub Actions, Terraform, Git, Kafka, R, Scala', '2026-02-24', 'https://www.indeed.com/viewjob?jk=16af568f8829f63e')
    ,('Senior Engineer I, Observability', 'DigitalOcean', 'Boston, MA, US USA', 10, 'RAG, Kubernetes, CI/CD, GitHub Actions, Terraform, Git, Kafka, R, Scala', '2026-02-24', 'https://www.indeed.com/viewjob?jk=1dd9190af8f5e921')
    ,('Senior Engineer I, Observability', 'DigitalOcean', 'Denver, CO, US USA', 10, 'RAG, Kubernetes, CI/CD, GitHub Actions, Terraform, Git, Kafka, R, Scala', '2026-02-24', 'https://www.indeed.com/viewjob?jk=faf1fcd9ebafc0c9')
    ,('Senior Engineer I, Observability', 'DigitalOcean', 'Austin, TX, US USA', 10, 'RAG, Kubernetes, CI/CD, GitHub Actions, Terraform, Git, Kafka, R, Scala', '2026-02-24', 'https://www.indeed.com/viewjob?jk=ec3420072b142500')
    ,('Senior Engineer I, Observability', 'DigitalOcean', 'San Francisco, CA, US USA', 10, 'RAG, Kubernetes, CI/CD, GitHub Actions, Terraform, Git, Kafka, R, Scala', '2026-02-24', 'https://www.indeed.com/viewjob?jk=93d6453a656d5be2')
    ,('Engineer III (Software EUP)', 'Ross Dress For Less', 'Dublin, CA, US USA', 10, 'RAG, Copilot, CI/CD, Jenkins, GitHub Actions, Git, R, Scala, Optimization', '2026-02-24', 'https://www.indeed.com/viewjob?jk=38362cff165251d5')
    ,('AI/GIS Engineer', 'Alson Services', 'Grand Rapids, MI, US USA', 10, 'LangChain, RAG, PostgreSQL, MongoDB, Python, SQL, R, Java, Scala', '2026-02-24', 'https://www.indeed.com/viewjob?jk=6b2577bc005b4c8f')
    ,('Software Engineer II, Data Services', 'Fox Corporation', 'Los Angeles, CA, US USA', 10, 'RAG, Kubernetes, CI/CD, GitHub Actions, Terraform, Git, Kafka, Python, R', '2026-02-24', 'https://www.indeed.com/viewjob?jk=717fdee43693f735')
    ,('AI Engineer I', 'BillGO, Inc.', 'Fort Collins, CO, US USA', 10, 'AI Engineer, Prompt Engineering, CI/CD, NoSQL, Python, SQL, R, Java, Scala', '2026-02-24', 'https://www.indeed.com/viewjob?jk=72139005356b854a')
    ,('Data Scientist', 'Qode', 'Basking Ridge, NJ, US USA', 10, 'Data Scientist, XGBoost, BigQuery, BigQuery, Tableau, Python, SQL, R, Optimization', '2026-02-24', 'https://www.indeed.com/viewjob?jk=31fc8ee741d634cf')
    ,('Senior Data Scientist, Condition Monitoring', 'Caterpillar', 'Chicago, IL, US USA', 10, 'Data Scientist, RAG, PyTorch, Git, Snowflake, Python, R, Optimization, Bayesian', '2026-02-24', 'https://www.indeed.com/viewjob?jk=fc8e43d1baa05943')
)

# Ensure the Posted At column (F) stays plain text, not auto-converted to a date
$ws.Range("F2:F24").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
